# Tidsregistrering i PTE projektet - Nada Omer
# Apply the changes described in the commit:
# "Har tilføjet implements til NormalspaendingImpl - importeret kalsser til OC5Test"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tidsregistrering")

# --- Update a handful of existing activity / duration cells (rows 33, 36-40) ---
$ws.Range("F33").Value = "Implementer OC4"

$ws.Range("F36").Value = "Krydstjeck for UC7"
$ws.Range("F37").Value = "Krydstjeck for UC8"
$ws.Range("F38").Value = "Review kode til OC3"
$ws.Range("F39").Value = "Review kode til OC1"

$ws.Range("F40").Value = "Review kode til OC2"
$ws.Range("I40").Value = "0 time : 15 min."

# --- Add new time-registration rows 41-49, copying the formatting from row 40 ---
$ws.Range("A40:C40").Copy()
$ws.Range("A41:C49").PasteSpecial(-4122)
$ws.Range("E40:I40").Copy()
$ws.Range("E41:I48").PasteSpecial(-4122)
# Row 49 only uses columns A-G (no Sluttid/Samlet tid entries yet)
$ws.Range("E40:G40").Copy()
$ws.Range("E49:G49").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row 41
$ws.Range("A41").Value = 42809
$ws.Range("B41").Value = "NO"
$ws.Range("C41").Value = "Nada H. A. Omer"
$ws.Range("E41").Value = "Reviewer"
$ws.Range("F41").Value = "Review OC7"
$ws.Range("G41").Value = 0.35069444444444442
$ws.Range("H41").Value = 0.375
$ws.Range("I41").Value = "0 time : 35 min."

# Row 42
$ws.Range("A42").Value = 42809
$ws.Range("B42").Value = "NO"
$ws.Range("C42").Value = "Nada H. A. Omer"
$ws.Range("E42").Value = "Designer"
$ws.Range("F42").Value = "SD og DCD for OC7"
$ws.Range("G42").Value = 0.3888888888888889
$ws.Range("H42").Value = 0.44444444444444442
$ws.Range("I42").Value = "1 time : 20 min."

# Row 43
$ws.Range("A43").Value = 42809
$ws.Range("B43").Value = "NO"
$ws.Range("C43").Value = "Nada H. A. Omer"
$ws.Range("E43").Value = "Implementer"
$ws.Range("F43").Value = "Implementer Junit test til OC5"
$ws.Range("G43").Value = 0.44444444444444442
$ws.Range("H43").Value = 0.52777777777777779
$ws.Range("I43").Value = "1 time : 30 min."

# Row 44
$ws.Range("A44").Value = 42809
$ws.Range("B44").Value = "NO"
$ws.Range("C44").Value = "Nada H. A. Omer"
$ws.Range("E44").Value = "Implementer"
$ws.Range("F44").Value = "Prøver at implemente Junit test til OC7, men Testsutie var indvalid."
$ws.Range("G44").Value = 0.53125
$ws.Range("H44").Value = 0.57291666666666663
$ws.Range("I44").Value = "1 time : 00 min."

# Row 45
$ws.Range("A45").Value = 42809
$ws.Range("B45").Value = "NO"
$ws.Range("C45").Value = "Nada H. A. Omer"
$ws.Range("E45").Value = "Implementer"
$ws.Range("F45").Value = "Prøver at implemente OC5, men Design var indvalid."
$ws.Range("G45").Value = 0.57291666666666663
$ws.Range("H45").Value = 0.59722222222222221
$ws.Range("I45").Value = "0 time : 35 min."

# Row 46
$ws.Range("A46").Value = 42809
$ws.Range("B46").Value = "NO"
$ws.Range("C46").Value = "Nada H. A. Omer"
$ws.Range("E46").Value = "Implementer"
$ws.Range("F46").Value = "Implementer OC7"
$ws.Range("G46").Value = 0.60069444444444442
$ws.Range("H46").Value = 0.63541666666666663
$ws.Range("I46").Value = "0 time : 50 min."

# Row 47
$ws.Range("A47").Value = 42810
$ws.Range("B47").Value = "NO"
$ws.Range("C47").Value = "Nada H. A. Omer"
$ws.Range("E47").Value = "Implementer"
$ws.Range("F47").Value = "Implementer OC7"
$ws.Range("G47").Value = 0.33680555555555558
$ws.Range("H47").Value = 0.38541666666666669
$ws.Range("I47").Value = "1 time : 10 min."

# Row 48
$ws.Range("A48").Value = 42810
$ws.Range("B48").Value = "NO"
$ws.Range("C48").Value = "Nada H. A. Omer"
$ws.Range("E48").Value = "Implementer"
$ws.Range("F48").Value = "Implementer OC5"
$ws.Range("G48").Value = 0.3888888888888889
$ws.Range("H48").Value = 0.41666666666666669
$ws.Range("I48").Value = "0 time : 40 min."

# Row 49 (shorter row - only through column G)
$ws.Range("A49").Value = 42810
$ws.Range("B49").Value = "NO"
$ws.Range("C49").Value = "Nada H. A. Omer"
$ws.Range("E49").Value = "Implementer"
$ws.Range("F49").Value = "Implementer  Junit test til OC7"
$ws.Range("G49").Value = 0.4236111111111111

# --- Column F got wider to fit the longer activity descriptions ---
$ws.Columns.Item(6).ColumnWidth = 58.6

# --- Data validation lists: re-point at the now-larger ranges, Deltagere first ---
$ws.Range("E3:E122").Validation.Delete()
$ws.Range("C3:C41").Validation.Delete()
$ws.Range("C3:C49").Validation.Add(3, 1, 1, "=Deltagere")
$ws.Range("E3:E124").Validation.Add(3, 1, 1, "=GyldigeRoller")

# --- Move the selection to reflect the new bottom of the table ---
$ws.Range("I50").Select()
